{"js": "const replacements = [\n  [\"93\u00f74=\", \"75\u00f75=\"],\n  [\"86\u00f75=\", \"99\u00f78=\"],\n  [\"95\u00f73=\", \"25\u00f73=\"],\n  [\"78\u00f74=\", \"71\u00f73=\"],\n  [\"51\u00f79=\", \"60\u00f72=\"],\n  [\"97\u00f76=\", \"73\u00f75=\"],\n  [\"18\u00f78=\", \"68\u00f74=\"],\n  [\"85\u00f76=\", \"47\u00f77=\"],\n  [\"53\u00f78=\", \"54\u00f75=\"],\n  [\"64\u00f79=\", \"22\u00f77=\"],\n  [\"57\u00f76=\", \"40\u00f77=\"],\n  [\"80\u00f78=\", \"82\u00f77=\"],\n  [\"29\u00f72=\", \"71\u00f79=\"],\n  [\"97\u00f78=\", \"59\u00f75=\"],\n  [\"37\u00f78=\", \"96\u00f75=\"],\n  [\"90\u00f74=\", \"53\u00f78=\"],\n  [\"50\u00f76=\", \"36\u00f76=\"],\n  [\"33\u00f73=\", \"33\u00f77=\"],\n  [\"76\u00f74=\", \"31\u00f72=\"],\n  [\"53\u00f73=\", \"59\u00f72=\"],\n  [\"80\u00f72=\", \"31\u00f78=\"],\n  [\"55\u00f77=\", \"15\u00f75=\"],\n  [\"95\u00f75=\", \"76\u00f76=\"],\n  [\"33\u00f74=\", \"23\u00f77=\"],\n  [\"98\u00f75=\", \"59\u00f78=\"],\n];\n\nfor (const [before, after] of replacements) {\n  const results = context.document.body.search(before, { matchCase: true, matchWholeWord: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly 1 match for '\" + before + \"' but found \" + results.items.length);\n  }\n  results.items[0].insertText(after, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"93\u00f74=\", \"75\u00f75=\"),\n    @(\"86\u00f75=\", \"99\u00f78=\"),\n    @(\"95\u00f73=\", \"25\u00f73=\"),\n    @(\"78\u00f74=\", \"71\u00f73=\"),\n    @(\"51\u00f79=\", \"60\u00f72=\"),\n    @(\"97\u00f76=\", \"73\u00f75=\"),\n    @(\"18\u00f78=\", \"68\u00f74=\"),\n    @(\"85\u00f76=\", \"47\u00f77=\"),\n    @(\"53\u00f78=\", \"54\u00f75=\"),\n    @(\"64\u00f79=\", \"22\u00f77=\"),\n    @(\"57\u00f76=\", \"40\u00f77=\"),\n    @(\"80\u00f78=\", \"82\u00f77=\"),\n    @(\"29\u00f72=\", \"71\u00f79=\"),\n    @(\"97\u00f78=\", \"59\u00f75=\"),\n    @(\"37\u00f78=\", \"96\u00f75=\"),\n    @(\"90\u00f74=\", \"53\u00f78=\"),\n    @(\"50\u00f76=\", \"36\u00f76=\"),\n    @(\"33\u00f73=\", \"33\u00f77=\"),\n    @(\"76\u00f74=\", \"31\u00f72=\"),\n    @(\"53\u00f73=\", \"59\u00f72=\"),\n    @(\"80\u00f72=\", \"31\u00f78=\"),\n    @(\"55\u00f77=\", \"15\u00f75=\"),\n    @(\"95\u00f75=\", \"76\u00f76=\"),\n    @(\"33\u00f74=\", \"23\u00f77=\"),\n    @(\"98\u00f75=\", \"59\u00f78=\"),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    # MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace(=wdReplaceAll)\n    $find.Execute($pair[0], $true, $true, $false, $false, $false, $true, $null, $false, $pair[1], 2)\n}\n"}
